# Consolidate adjacent plain-text runs in slide titles into single runs,
# matching the PowerPoint writer's run-merging behavior described in the
# commit message. Each title's trailing formatted run (e.g. "inline code"
# in Consolas) is left untouched; only the leading plain-text runs that
# were unnecessarily split are merged via a Characters() sub-range set.

$p = $ppt.ActivePresentation

# Slide 1 title: "Header" + " " + "with" + " " -> "Header with "
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(1, 12).Text = "Header with "

# Slide 2 title: "Syntax" + " " + "highlighting" -> "Syntax highlighting"
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2.Characters(1, $tr2.Length).Text = "Syntax highlighting"

# Slide 3 title: "Two" + " " + "column" + " " + "slide" -> "Two column slide"
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Characters(1, $tr3.Length).Text = "Two column slide"
